$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").Value = '26.878.76'
$ws.Range("E2").Value = '  +0.10%  '
$ws.Range("D3").Value = '1.640.12'
$ws.Range("E4").Value = '  -0.55%  '
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '216.89'
$c.Style = "Normal"
$ws.Range("E5").Value = '  -0.79%  '
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '0.513'
$c.Style = "Normal"
$ws.Range("E6").Value = '  +2.52%  '
$ws.Range("E7").Value = '  -0.55%  '
$ws.Range("E8").Value = '  +1.75%  '
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '0.0626'
$c.Style = "Normal"
$ws.Range("E9").Value = '  +0.42%  '
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '19.90'
$c.Style = "Normal"
$ws.Range("E10").Value = '  +3.31%  '
$ws.Range("D12").Value = '1.869.43'
$ws.Range("E12").Value = '  -0.17%  '
$ws.Range("D13").Value = '1.634.58'
$ws.Range("E13").Value = '  -0.59%  '
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '4.13'
$c.Style = "Normal"
$ws.Range("E14").Value = '  -0.78%  '
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '0.531'
$c.Style = "Normal"
$ws.Range("E15").Value = '  +0.89%  '
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '67.47'
$c.Style = "Normal"
$ws.Range("E16").Value = '  +3.16%  '
$ws.Range("D17").Value = '26.868.15'
$ws.Range("D18").Value = '0.0₃0732'
$ws.Range("E18").Value = '  -0.49%  '
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '219.66'
$c.Style = "Normal"
$ws.Range("E19").Value = '  +1.69%  '
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '6.85'
$c.Style = "Normal"
$ws.Range("E21").Value = '  +3.18%  '
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '4.40'
$c.Style = "Normal"
$ws.Range("E22").Value = '  +0.59%  '
$ws.Range("E23").Value = '  +3.27%  '
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '9.16'
$c.Style = "Normal"
$ws.Range("E24").Value = '  -0.57%  '
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '147.36'
$c.Style = "Normal"
$ws.Range("E26").Value = '  -0.60%  '
$ws.Range("E27").Value = '  +2.73%  '
$ws.Range("E28").Value = '  +0.33%  '
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '15.82'
$c.Style = "Normal"
$ws.Range("E29").Value = '  +0.39%  '
$ws.Range("E30").Value = '  -1.31%  '
$ws.Range("E31").Value = '  -0.94%  '
$ws.Range("E32").Value = '  -1.47%  '
$ws.Range("E33").Value = '  +0.33%  '
$ws.Range("E34").Value = '  +1.40%  '
$ws.Range("D35").Value = '1.270.10'
$ws.Range("E35").Value = '  -0.17%  '
$ws.Range("E36").Value = '  -0.08%  '
$ws.Range("E37").Value = '  +1.80%  '
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '0.535'
$c.Style = "Normal"
$ws.Range("E38").Value = '  +0.38%  '
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '0.834'
$c.Style = "Normal"
$ws.Range("E39").Value = '  +2.17%  '
$ws.Range("E40").Value = '  -0.54%  '
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '0.809'
$c.Style = "Normal"
$ws.Range("E41").Value = '  +0.63%  '
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '5.39'
$c.Style = "Normal"
$ws.Range("E42").Value = '  +0.70%  '
$ws.Range("D43").Value = '1.780.72'
$ws.Range("E43").Value = '  -0.15%  '
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '62.02'
$c.Style = "Normal"
$ws.Range("E44").Value = '  +1.23%  '
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '2.10'
$c.Style = "Normal"
$ws.Range("E45").Value = '  -1.42%  '
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '91.71'
$c.Style = "Normal"
$ws.Range("E46").Value = '  -1.23%  '
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '1.59'
$c.Style = "Normal"
$ws.Range("E47").Value = '  -1.19%  '
$ws.Range("D48").Value = '0.0⁦0105'
$ws.Range("E48").Value = '  +2.79%  '
$ws.Range("E49").Value = '  -0.68%  '
$ws.Range("E50").Value = '  +1.31%  '
$ws.Range("E51").Value = '  -0.29%  '
